$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = -22.192
$ws.Range("C4").Value = -12.605
$ws.Range("D4").Value = -7.703

# Row 5
$ws.Range("C5").Value = -12.413

# Row 7
$ws.Range("A7").Value = -20.987

# Row 8
$ws.Range("C8").Value = -12.8

# Row 9
$ws.Range("D9").Value = -7.867999999999999

# Row 16
$ws.Range("A16").Value = -21.148
$ws.Range("C16").Value = -12.759

# Row 18
$ws.Range("D18").Value = -7.835999999999999
